# Loan RBI, Variable Instalments
#
# On the "Repayment schedule" sheet, a new (blank) column is inserted
# before the existing "Late" column (column N), pushing the old N/O/P
# columns (Late / Outstanding-heading / Outstanding) one column to the
# right (O/P/Q). The new column takes on the same width as the column
# to its left ("In Advance", column M). Finally, the "Repayment
# schedule" sheet (not "Transactions") is left as the active sheet /
# selected tab, with a trailing selection a bit to the right of the
# data that was just inserted.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N ("Late"), shifting the
# existing N/O/P columns (and all their data/styles) one place to the
# right.
$ws.Columns("N").Insert() | Out-Null

# Give the newly inserted column the same width as column M
# ("In Advance"), matching the width Excel would apply automatically.
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Make "Repayment schedule" the active sheet (was "Transactions"
# before), and leave the selection on the cell the user clicked after
# inserting the column.
$ws.Activate() | Out-Null
$ws.Range("R8").Select() | Out-Null
